$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Device" column (F) of the Photo Report table.
# Rows 18 and 19 are special-cased: their device is changed to "Samsung".
# All other rows whose device is currently "LG", "HTC", or "Motorola" are
# changed to the new device model "Tecno".
for ($r = 2; $r -le 23; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $val = $cell.Value2

    if ($r -eq 18 -or $r -eq 19) {
        if ($val -eq "LG" -or $val -eq "HTC" -or $val -eq "Motorola") {
            $cell.Value = "Samsung"
        }
    }
    elseif ($val -eq "LG" -or $val -eq "HTC" -or $val -eq "Motorola") {
        $cell.Value = "Tecno"
    }
}

$ws.Range("I18").Select() | Out-Null
